$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the product_type for "Tenra Basho: Vampire Princess" (row 5) from "scenario" to "box set"
$ws.Range("F5").Value = "box set"

# Update the active cell selection shown in the sheet view
$ws.Range("F6").Select()
